# fixes in upload count/close
# Adds 3 new regression-run rows (52-54) to the AMSIN sheet and fixes up
# the formatting / run-time value of the existing last row (51).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "AMSIN" sheet

# ---------------------------------------------------------------------
# Step 1: snapshot the (currently unstyled) row 51 formatting onto the
# new row 54, before row 51 itself gets restyled below. In the original
# file row 51 is the only data row whose cells carry no explicit style
# (A/C/D/E/F/G have no "s" attribute, B51 uses s="12"); row 54 needs
# that exact same "no explicit style" look, so we copy it now while it
# still exists.
# ---------------------------------------------------------------------
$ws.Range("A51:G51").Copy()
$ws.Range("A54:G54").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------
# Step 2: fix up existing row 51 - give it the normal data-row styling
# (same as row 50, style index 6) and correct B51's run-time value.
# ---------------------------------------------------------------------
$ws.Range("A51").Clear()
$ws.Range("C51:G51").Clear()

$ws.Range("A50").Copy()
$ws.Range("A51").PasteSpecial(-4163)       # xlPasteValues (keeps A51 as text "2022-02-24")

$ws.Range("B51").Value = 44616.55075899306
$ws.Range("C51").Value = "asa2332"
$ws.Range("D51").Value = 269
$ws.Range("E51").Value = 0
$ws.Range("F51").Value = 269
$ws.Range("G51").Value = 0.26

# ---------------------------------------------------------------------
# Step 3: new row 52
# ---------------------------------------------------------------------
$ws.Range("A50").Copy()
$ws.Range("A52").PasteSpecial(-4163)       # text "2022-02-24"

$ws.Range("B51").Copy()
$ws.Range("B52").PasteSpecial(-4122)       # style 12 (date/time number format)
$ws.Range("B52").Value = 44616.70045383102

$ws.Range("C52").Value = "test1234"
$ws.Range("D52").Value = 269
$ws.Range("E52").Value = 0
$ws.Range("F52").Value = 269
$ws.Range("G52").Value = 0.26

# ---------------------------------------------------------------------
# Step 4: new row 53 - "2022-02-25" does not exist anywhere yet, so
# build it via a text formula and flatten the formula to a literal
# value with a self copy/paste-values (avoids Excel's auto date
# conversion that a direct .Value assignment of a date-shaped string
# would trigger).
# ---------------------------------------------------------------------
$ws.Range("A53").Formula = "=""2022-02-25"""
$ws.Range("A53").Copy()
$ws.Range("A53").PasteSpecial(-4163)

$ws.Range("B51").Copy()
$ws.Range("B53").PasteSpecial(-4122)
$ws.Range("B53").Value = 44617.49754053241

$ws.Range("C53").Value = "test158"
$ws.Range("D53").Value = 269
$ws.Range("E53").Value = 267
$ws.Range("F53").Value = 2
$ws.Range("G53").Value = 5.45

# ---------------------------------------------------------------------
# Step 5: row 54 - formatting was already applied in Step 1, now only
# fill in the values (same text-flattening trick for the date).
# ---------------------------------------------------------------------
$ws.Range("A54").Formula = "=""2022-02-25"""
$ws.Range("A54").Copy()
$ws.Range("A54").PasteSpecial(-4163)

$ws.Range("B54").Value = 44617.52886187014

$ws.Range("C54").Value = "test158"
$ws.Range("D54").Value = 269
$ws.Range("E54").Value = 74
$ws.Range("F54").Value = 195
$ws.Range("G54").Value = 3.21
